$d = $word.ActiveDocument

$d.Content.Find.Execute("25+62=87", $true, $false, $false, $false, $false, $true, 1, $false, "47+39=86", 2) | Out-Null
$d.Content.Find.Execute("47+6=53", $true, $false, $false, $false, $false, $true, 1, $false, "1+90=91", 2) | Out-Null
$d.Content.Find.Execute("70-18=52", $true, $false, $false, $false, $false, $true, 1, $false, "97-7=90", 2) | Out-Null
$d.Content.Find.Execute("30+0=30", $true, $false, $false, $false, $false, $true, 1, $false, "25-4=21", 2) | Out-Null
$d.Content.Find.Execute("80-44=36", $true, $false, $false, $false, $false, $true, 1, $false, "90-56=34", 2) | Out-Null
$d.Content.Find.Execute("89-16=73", $true, $false, $false, $false, $false, $true, 1, $false, "83-48=35", 2) | Out-Null
$d.Content.Find.Execute("41+48=89", $true, $false, $false, $false, $false, $true, 1, $false, "5+35=40", 2) | Out-Null
$d.Content.Find.Execute("27+14=41", $true, $false, $false, $false, $false, $true, 1, $false, "13+38=51", 2) | Out-Null
$d.Content.Find.Execute("58-11=47", $true, $false, $false, $false, $false, $true, 1, $false, "47-10=37", 2) | Out-Null
$d.Content.Find.Execute("46-21=25", $true, $false, $false, $false, $false, $true, 1, $false, "96-83=13", 2) | Out-Null
$d.Content.Find.Execute("49-5=44", $true, $false, $false, $false, $false, $true, 1, $false, "1+66=67", 2) | Out-Null
$d.Content.Find.Execute("95-28=67", $true, $false, $false, $false, $false, $true, 1, $false, "94+1=95", 2) | Out-Null
$d.Content.Find.Execute("62-57=5", $true, $false, $false, $false, $false, $true, 1, $false, "51+24=75", 2) | Out-Null
$d.Content.Find.Execute("40-24=16", $true, $false, $false, $false, $false, $true, 1, $false, "5+68=73", 2) | Out-Null
$d.Content.Find.Execute("0+41=41", $true, $false, $false, $false, $false, $true, 1, $false, "18+45=63", 2) | Out-Null
$d.Content.Find.Execute("37+19=56", $true, $false, $false, $false, $false, $true, 1, $false, "78-3=75", 2) | Out-Null
$d.Content.Find.Execute("45+49=94", $true, $false, $false, $false, $false, $true, 1, $false, "60+11=71", 2) | Out-Null
$d.Content.Find.Execute("33-10=23", $true, $false, $false, $false, $false, $true, 1, $false, "3+12=15", 2) | Out-Null
$d.Content.Find.Execute("12+85=97", $true, $false, $false, $false, $false, $true, 1, $false, "65-39=26", 2) | Out-Null
$d.Content.Find.Execute("87-17=70", $true, $false, $false, $false, $false, $true, 1, $false, "27+30=57", 2) | Out-Null
$d.Content.Find.Execute("94-94=0", $true, $false, $false, $false, $false, $true, 1, $false, "68-45=23", 2) | Out-Null
$d.Content.Find.Execute("21-9=12", $true, $false, $false, $false, $false, $true, 1, $false, "9-5=4", 2) | Out-Null
$d.Content.Find.Execute("4+86=90", $true, $false, $false, $false, $false, $true, 1, $false, "47-19=28", 2) | Out-Null
$d.Content.Find.Execute("40+53=93", $true, $false, $false, $false, $false, $true, 1, $false, "38-3=35", 2) | Out-Null
$d.Content.Find.Execute("1+44=45", $true, $false, $false, $false, $false, $true, 1, $false, "52-49=3", 2) | Out-Null
$d.Content.Find.Execute("91-85=6", $true, $false, $false, $false, $false, $true, 1, $false, "24+12=36", 2) | Out-Null
$d.Content.Find.Execute("13-6=7", $true, $false, $false, $false, $false, $true, 1, $false, "82-21=61", 2) | Out-Null
$d.Content.Find.Execute("69-7=62", $true, $false, $false, $false, $false, $true, 1, $false, "55-53=2", 2) | Out-Null
$d.Content.Find.Execute("90-12=78", $true, $false, $false, $false, $false, $true, 1, $false, "9-1=8", 2) | Out-Null
$d.Content.Find.Execute("30+55=85", $true, $false, $false, $false, $false, $true, 1, $false, "89-48=41", 2) | Out-Null
$d.Content.Find.Execute("29+27=56", $true, $false, $false, $false, $false, $true, 1, $false, "20+53=73", 2) | Out-Null
$d.Content.Find.Execute("53+13=66", $true, $false, $false, $false, $false, $true, 1, $false, "10+30=40", 2) | Out-Null
$d.Content.Find.Execute("88-54=34", $true, $false, $false, $false, $false, $true, 1, $false, "72-45=27", 2) | Out-Null
$d.Content.Find.Execute("58+25=83", $true, $false, $false, $false, $false, $true, 1, $false, "11+68=79", 2) | Out-Null
$d.Content.Find.Execute("82-77=5", $true, $false, $false, $false, $false, $true, 1, $false, "75-12=63", 2) | Out-Null
$d.Content.Find.Execute("18+78=96", $true, $false, $false, $false, $false, $true, 1, $false, "97-21=76", 2) | Out-Null
$d.Content.Find.Execute("22-0=22", $true, $false, $false, $false, $false, $true, 1, $false, "50+4=54", 2) | Out-Null
$d.Content.Find.Execute("21+3=24", $true, $false, $false, $false, $false, $true, 1, $false, "26+43=69", 2) | Out-Null
$d.Content.Find.Execute("1+50=51", $true, $false, $false, $false, $false, $true, 1, $false, "46+2=48", 2) | Out-Null
$d.Content.Find.Execute("56-46=10", $true, $false, $false, $false, $false, $true, 1, $false, "86+13=99", 2) | Out-Null
$d.Content.Find.Execute("95-61=34", $true, $false, $false, $false, $false, $true, 1, $false, "15-8=7", 2) | Out-Null
$d.Content.Find.Execute("29-29=0", $true, $false, $false, $false, $false, $true, 1, $false, "90+5=95", 2) | Out-Null
$d.Content.Find.Execute("90-41=49", $true, $false, $false, $false, $false, $true, 1, $false, "26+58=84", 2) | Out-Null
$d.Content.Find.Execute("69+20=89", $true, $false, $false, $false, $false, $true, 1, $false, "56-22=34", 2) | Out-Null
$d.Content.Find.Execute("10+14=24", $true, $false, $false, $false, $false, $true, 1, $false, "26+69=95", 2) | Out-Null
$d.Content.Find.Execute("89-20=69", $true, $false, $false, $false, $false, $true, 1, $false, "48-5=43", 2) | Out-Null
$d.Content.Find.Execute("94-15=79", $true, $false, $false, $false, $false, $true, 1, $false, "77-71=6", 2) | Out-Null
$d.Content.Find.Execute("38-27=11", $true, $false, $false, $false, $false, $true, 1, $false, "79-26=53", 2) | Out-Null
$d.Content.Find.Execute("69-68=1", $true, $false, $false, $false, $false, $true, 1, $false, "52+39=91", 2) | Out-Null
$d.Content.Find.Execute("59-22=37", $true, $false, $false, $false, $false, $true, 1, $false, "20+20=40", 2) | Out-Null
$d.Content.Find.Execute("2+19=21", $true, $false, $false, $false, $false, $true, 1, $false, "24+30=54", 2) | Out-Null
$d.Content.Find.Execute("90-79=11", $true, $false, $false, $false, $false, $true, 1, $false, "28-2=26", 2) | Out-Null
$d.Content.Find.Execute("35+3=38", $true, $false, $false, $false, $false, $true, 1, $false, "7+76=83", 2) | Out-Null
$d.Content.Find.Execute("79+10=89", $true, $false, $false, $false, $false, $true, 1, $false, "27+40=67", 2) | Out-Null
$d.Content.Find.Execute("27+32=59", $true, $false, $false, $false, $false, $true, 1, $false, "8+40=48", 2) | Out-Null
$d.Content.Find.Execute("98-70=28", $true, $false, $false, $false, $false, $true, 1, $false, "25-15=10", 2) | Out-Null
$d.Content.Find.Execute("56+15=71", $true, $false, $false, $false, $false, $true, 1, $false, "4+21=25", 2) | Out-Null
$d.Content.Find.Execute("91-50=41", $true, $false, $false, $false, $false, $true, 1, $false, "75+4=79", 2) | Out-Null
$d.Content.Find.Execute("48+25=73", $true, $false, $false, $false, $false, $true, 1, $false, "86-21=65", 2) | Out-Null
$d.Content.Find.Execute("89-31=58", $true, $false, $false, $false, $false, $true, 1, $false, "70-37=33", 2) | Out-Null
$d.Content.Find.Execute("18+46=64", $true, $false, $false, $false, $false, $true, 1, $false, "51-1=50", 2) | Out-Null
$d.Content.Find.Execute("66+29=95", $true, $false, $false, $false, $false, $true, 1, $false, "53+35=88", 2) | Out-Null
$d.Content.Find.Execute("12+78=90", $true, $false, $false, $false, $false, $true, 1, $false, "75-52=23", 2) | Out-Null
$d.Content.Find.Execute("24+51=75", $true, $false, $false, $false, $false, $true, 1, $false, "10+6=16", 2) | Out-Null
$d.Content.Find.Execute("12+29=41", $true, $false, $false, $false, $false, $true, 1, $false, "23+25=48", 2) | Out-Null
$d.Content.Find.Execute("24-9=15", $true, $false, $false, $false, $false, $true, 1, $false, "20+25=45", 2) | Out-Null
$d.Content.Find.Execute("45+37=82", $true, $false, $false, $false, $false, $true, 1, $false, "89-32=57", 2) | Out-Null
$d.Content.Find.Execute("12+37=49", $true, $false, $false, $false, $false, $true, 1, $false, "8+69=77", 2) | Out-Null
$d.Content.Find.Execute("17-1=16", $true, $false, $false, $false, $false, $true, 1, $false, "32+35=67", 2) | Out-Null
$d.Content.Find.Execute("16+21=37", $true, $false, $false, $false, $false, $true, 1, $false, "63+23=86", 2) | Out-Null
$d.Content.Find.Execute("15+10=25", $true, $false, $false, $false, $false, $true, 1, $false, "13+12=25", 2) | Out-Null
$d.Content.Find.Execute("59-42=17", $true, $false, $false, $false, $false, $true, 1, $false, "8+77=85", 2) | Out-Null
$d.Content.Find.Execute("17-0=17", $true, $false, $false, $false, $false, $true, 1, $false, "55+32=87", 2) | Out-Null
$d.Content.Find.Execute("93-20=73", $true, $false, $false, $false, $false, $true, 1, $false, "89-46=43", 2) | Out-Null
$d.Content.Find.Execute("68-7=61", $true, $false, $false, $false, $false, $true, 1, $false, "84-63=21", 2) | Out-Null
$d.Content.Find.Execute("48-9=39", $true, $false, $false, $false, $false, $true, 1, $false, "0+35=35", 2) | Out-Null
$d.Content.Find.Execute("62-50=12", $true, $false, $false, $false, $false, $true, 1, $false, "39+26=65", 2) | Out-Null
$d.Content.Find.Execute("31+36=67", $true, $false, $false, $false, $false, $true, 1, $false, "84-26=58", 2) | Out-Null
$d.Content.Find.Execute("27+59=86", $true, $false, $false, $false, $false, $true, 1, $false, "8+75=83", 2) | Out-Null
$d.Content.Find.Execute("11+21=32", $true, $false, $false, $false, $false, $true, 1, $false, "3+80=83", 2) | Out-Null
$d.Content.Find.Execute("49-32=17", $true, $false, $false, $false, $false, $true, 1, $false, "2+65=67", 2) | Out-Null
$d.Content.Find.Execute("48+43=91", $true, $false, $false, $false, $false, $true, 1, $false, "3+83=86", 2) | Out-Null
$d.Content.Find.Execute("88-65=23", $true, $false, $false, $false, $false, $true, 1, $false, "19-19=0", 2) | Out-Null
$d.Content.Find.Execute("7+36=43", $true, $false, $false, $false, $false, $true, 1, $false, "60-26=34", 2) | Out-Null
$d.Content.Find.Execute("88-80=8", $true, $false, $false, $false, $false, $true, 1, $false, "85-44=41", 2) | Out-Null
$d.Content.Find.Execute("84-78=6", $true, $false, $false, $false, $false, $true, 1, $false, "28+21=49", 2) | Out-Null
$d.Content.Find.Execute("0+14=14", $true, $false, $false, $false, $false, $true, 1, $false, "17+38=55", 2) | Out-Null
$d.Content.Find.Execute("15+8=23", $true, $false, $false, $false, $false, $true, 1, $false, "66-12=54", 2) | Out-Null
$d.Content.Find.Execute("39+46=85", $true, $false, $false, $false, $false, $true, 1, $false, "99-53=46", 2) | Out-Null
$d.Content.Find.Execute("47+2=49", $true, $false, $false, $false, $false, $true, 1, $false, "47+31=78", 2) | Out-Null
$d.Content.Find.Execute("28-0=28", $true, $false, $false, $false, $false, $true, 1, $false, "99-99=0", 2) | Out-Null
$d.Content.Find.Execute("42-41=1", $true, $false, $false, $false, $false, $true, 1, $false, "80-23=57", 2) | Out-Null
$d.Content.Find.Execute("76+14=90", $true, $false, $false, $false, $false, $true, 1, $false, "8+82=90", 2) | Out-Null
$d.Content.Find.Execute("94-0=94", $true, $false, $false, $false, $false, $true, 1, $false, "22+63=85", 2) | Out-Null
$d.Content.Find.Execute("63+16=79", $true, $false, $false, $false, $false, $true, 1, $false, "23-9=14", 2) | Out-Null
$d.Content.Find.Execute("81+3=84", $true, $false, $false, $false, $false, $true, 1, $false, "65-55=10", 2) | Out-Null
$d.Content.Find.Execute("13+82=95", $true, $false, $false, $false, $false, $true, 1, $false, "14-1=13", 2) | Out-Null
$d.Content.Find.Execute("79-65=14", $true, $false, $false, $false, $false, $true, 1, $false, "30+27=57", 2) | Out-Null
$d.Content.Find.Execute("10-9=1", $true, $false, $false, $false, $false, $true, 1, $false, "35-11=24", 2) | Out-Null
$d.Content.Find.Execute("94-6=88", $true, $false, $false, $false, $false, $true, 1, $false, "32+5=37", 2) | Out-Null
